$d = $word.ActiveDocument

# 1. Change the placeholder "Name of character" table cell text to "Orla"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Name of character", $true, $false, $false, $false, $false, $true, 1, $false, "Orla", 2)

# 2. Merge the split "affect" runs (with spell-check proofErr wrapping) back
#    into a single contiguous sentence.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Discuss how the player's gestures/interactivity has an affect on the visual experience. ", $true, $false, $false, $false, $false, $true, 1, $false, "Discuss how the player's gestures/interactivity has an affect on the visual experience. ", 2)
